# Scheduled market-data refresh for Sheets/Balmung_Profits.xlsx
# Updates the currentAveragePrice* / Leve/ItemPrice* / Profit columns (H:N)
# for the rows whose underlying Universalis market data changed since the
# last run. Values only -- no formulas are stored in these cells.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 5275.84
$ws.Cells.Item(69, 9).Value = 4999.8887
$ws.Cells.Item(69, 10).Value = 5431.0625
$ws.Cells.Item(69, 11).Value = 14999.6661
$ws.Cells.Item(69, 12).Value = 16293.1875
$ws.Cells.Item(69, 13).Value = -14125.6661
$ws.Cells.Item(69, 14).Value = -18041.1875

$ws.Cells.Item(72, 8).Value = 5275.84
$ws.Cells.Item(72, 9).Value = 4999.8887
$ws.Cells.Item(72, 10).Value = 5431.0625
$ws.Cells.Item(72, 11).Value = 44998.99830000001
$ws.Cells.Item(72, 12).Value = 48879.5625
$ws.Cells.Item(72, 13).Value = -40630.99830000001
$ws.Cells.Item(72, 14).Value = -57615.5625

$ws.Cells.Item(80, 8).Value = 37037496
$ws.Cells.Item(80, 9).Value = 71428860
$ws.Cells.Item(80, 11).Value = 214286580
$ws.Cells.Item(80, 13).Value = -214285582

$ws.Cells.Item(83, 8).Value = 37037496
$ws.Cells.Item(83, 9).Value = 71428860
$ws.Cells.Item(83, 11).Value = 642859740
$ws.Cells.Item(83, 13).Value = -642854748

$ws.Cells.Item(86, 8).Value = 58825644
$ws.Cells.Item(86, 9).Value = 76924760
$ws.Cells.Item(86, 11).Value = 76924760
$ws.Cells.Item(86, 13).Value = -76923637

$ws.Cells.Item(89, 8).Value = 58825644
$ws.Cells.Item(89, 9).Value = 76924760
$ws.Cells.Item(89, 11).Value = 384623800
$ws.Cells.Item(89, 13).Value = -384618184

$ws.Cells.Item(94, 8).Value = 1113.1666
$ws.Cells.Item(94, 9).Value = 1113.1666
$ws.Cells.Item(94, 11).Value = 1113.1666
$ws.Cells.Item(94, 13).Value = -662.1666

$ws.Cells.Item(134, 8).Value = 137329.17
$ws.Cells.Item(134, 10).Value = 146177.27
$ws.Cells.Item(134, 12).Value = 146177.27
$ws.Cells.Item(134, 14).Value = -156317.27

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2120366.2
$ws.Cells.Item(45, 9).Value = 2226334.5
$ws.Cells.Item(45, 10).Value = 1000
$ws.Cells.Item(45, 11).Value = 2226334.5
$ws.Cells.Item(45, 12).Value = 1000
$ws.Cells.Item(45, 13).Value = -2225957.5
$ws.Cells.Item(45, 14).Value = -1754

$ws.Cells.Item(61, 8).Value = 1283981.9
$ws.Cells.Item(61, 9).Value = 6125.85
$ws.Cells.Item(61, 10).Value = 2266948
$ws.Cells.Item(61, 11).Value = 6125.85
$ws.Cells.Item(61, 12).Value = 2266948
$ws.Cells.Item(61, 13).Value = -5913.85
$ws.Cells.Item(61, 14).Value = -2267372

$ws.Cells.Item(74, 8).Value = 896248.7
$ws.Cells.Item(74, 9).Value = 7459.5557
$ws.Cells.Item(74, 10).Value = 1396192.5
$ws.Cells.Item(74, 11).Value = 7459.5557
$ws.Cells.Item(74, 12).Value = 1396192.5
$ws.Cells.Item(74, 13).Value = -6585.5557
$ws.Cells.Item(74, 14).Value = -1397940.5

$ws.Cells.Item(77, 8).Value = 896248.7
$ws.Cells.Item(77, 9).Value = 7459.5557
$ws.Cells.Item(77, 10).Value = 1396192.5
$ws.Cells.Item(77, 11).Value = 37297.7785
$ws.Cells.Item(77, 12).Value = 6980962.5
$ws.Cells.Item(77, 13).Value = -32929.7785
$ws.Cells.Item(77, 14).Value = -6989698.5

$ws.Cells.Item(132, 8).Value = 3128656.5
$ws.Cells.Item(132, 9).Value = 2079.75
$ws.Cells.Item(132, 11).Value = 6239.25
$ws.Cells.Item(132, 13).Value = -3709.25

$ws.Cells.Item(136, 8).Value = 1283981.9
$ws.Cells.Item(136, 9).Value = 6125.85
$ws.Cells.Item(136, 10).Value = 2266948
$ws.Cells.Item(136, 11).Value = 18377.55
$ws.Cells.Item(136, 12).Value = 6800844
$ws.Cells.Item(136, 13).Value = -15827.55
$ws.Cells.Item(136, 14).Value = -6805944

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 76923256
$ws.Cells.Item(80, 9).Value = 166666830
$ws.Cells.Item(80, 10).Value = 184.14285
$ws.Cells.Item(80, 11).Value = 166666830
$ws.Cells.Item(80, 12).Value = 184.14285
$ws.Cells.Item(80, 13).Value = -166665832
$ws.Cells.Item(80, 14).Value = -2180.14285

$ws.Cells.Item(83, 8).Value = 76923256
$ws.Cells.Item(83, 9).Value = 166666830
$ws.Cells.Item(83, 10).Value = 184.14285
$ws.Cells.Item(83, 11).Value = 833334150
$ws.Cells.Item(83, 12).Value = 920.71425
$ws.Cells.Item(83, 13).Value = -833329158
$ws.Cells.Item(83, 14).Value = -10904.71425

$ws.Cells.Item(134, 8).Value = 21954014
$ws.Cells.Item(134, 9).Value = 2129.0688
$ws.Cells.Item(134, 11).Value = 6387.2064
$ws.Cells.Item(134, 13).Value = -3852.2064

$ws.Cells.Item(141, 8).Value = 52756.547
$ws.Cells.Item(141, 10).Value = 52756.547
$ws.Cells.Item(141, 12).Value = 52756.547
$ws.Cells.Item(141, 14).Value = -63116.547

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(36, 8).Value = 3388.7778
$ws.Cells.Item(36, 9).Value = 4083.3333
$ws.Cells.Item(36, 11).Value = 4083.3333
$ws.Cells.Item(36, 13).Value = -3695.3333

$ws.Cells.Item(40, 8).Value = 3388.7778
$ws.Cells.Item(40, 9).Value = 4083.3333
$ws.Cells.Item(40, 11).Value = 4083.3333
$ws.Cells.Item(40, 13).Value = -3923.3333

$ws.Cells.Item(58, 8).Value = 2849.6191
$ws.Cells.Item(58, 9).Value = 2718.111
$ws.Cells.Item(58, 10).Value = 2948.25
$ws.Cells.Item(58, 11).Value = 2718.111
$ws.Cells.Item(58, 12).Value = 2948.25
$ws.Cells.Item(58, 13).Value = -2515.111
$ws.Cells.Item(58, 14).Value = -3354.25

$ws.Cells.Item(121, 8).Value = 300000
$ws.Cells.Item(121, 9).Value = 300000
$ws.Cells.Item(121, 10).Value = 0
$ws.Cells.Item(121, 11).Value = 300000
$ws.Cells.Item(121, 12).Value = 0
$ws.Cells.Item(121, 13).Value = -298690
$ws.Cells.Item(121, 14).ClearContents()

$ws.Cells.Item(136, 8).Value = 2849.6191
$ws.Cells.Item(136, 9).Value = 2718.111
$ws.Cells.Item(136, 10).Value = 2948.25
$ws.Cells.Item(136, 11).Value = 8154.333
$ws.Cells.Item(136, 12).Value = 8844.75
$ws.Cells.Item(136, 13).Value = -5604.333
$ws.Cells.Item(136, 14).Value = -13944.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(46, 8).Value = 15301024
$ws.Cells.Item(46, 10).Value = 7100001
$ws.Cells.Item(46, 12).Value = 21300003
$ws.Cells.Item(46, 14).Value = -21300185

$ws.Cells.Item(98, 8).Value = 783.1429000000001
$ws.Cells.Item(98, 10).Value = 996.25
$ws.Cells.Item(98, 12).Value = 2988.75
$ws.Cells.Item(98, 14).Value = -5984.75

$ws.Cells.Item(107, 8).Value = 60606388
$ws.Cells.Item(107, 9).Value = 295.22223
$ws.Cells.Item(107, 11).Value = 885.66669
$ws.Cells.Item(107, 13).Value = 1034.33331

$ws.Cells.Item(129, 8).Value = 19049962
$ws.Cells.Item(129, 9).Value = 1136.8
$ws.Cells.Item(129, 10).Value = 29632642
$ws.Cells.Item(129, 11).Value = 3410.4
$ws.Cells.Item(129, 12).Value = 88897926
$ws.Cells.Item(129, 13).Value = 1589.6
$ws.Cells.Item(129, 14).Value = -88907926

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(12, 8).Value = 2000
$ws.Cells.Item(12, 9).Value = 2000
$ws.Cells.Item(12, 11).Value = 2000
$ws.Cells.Item(12, 13).Value = -1860

$ws.Cells.Item(64, 8).Value = 58666.668
$ws.Cells.Item(64, 10).Value = 58666.668
$ws.Cells.Item(64, 12).Value = 58666.668
$ws.Cells.Item(64, 14).Value = -59162.668

$ws.Cells.Item(67, 8).Value = 58666.668
$ws.Cells.Item(67, 10).Value = 58666.668
$ws.Cells.Item(67, 12).Value = 58666.668
$ws.Cells.Item(67, 14).Value = -60382.668

$ws.Cells.Item(102, 8).Value = 33334414
$ws.Cells.Item(102, 9).Value = 35715344
$ws.Cells.Item(102, 11).Value = 35715344
$ws.Cells.Item(102, 13).Value = -35713722

$ws.Cells.Item(105, 8).Value = 70670.5
$ws.Cells.Item(105, 10).Value = 70670.5
$ws.Cells.Item(105, 12).Value = 70670.5
$ws.Cells.Item(105, 14).Value = -77658.5

$ws.Cells.Item(107, 8).Value = 111458.336
$ws.Cells.Item(107, 9).Value = 111458.336
$ws.Cells.Item(107, 11).Value = 111458.336
$ws.Cells.Item(107, 13).Value = -109538.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 8317.454
$ws.Cells.Item(7, 9).Value = 2566.8
$ws.Cells.Item(7, 11).Value = 2566.8
$ws.Cells.Item(7, 13).Value = -2454.8

$ws.Cells.Item(82, 8).Value = 3569.524
$ws.Cells.Item(82, 10).Value = 8164
$ws.Cells.Item(82, 12).Value = 8164
$ws.Cells.Item(82, 14).Value = -8886

$ws.Cells.Item(85, 8).Value = 3569.524
$ws.Cells.Item(85, 10).Value = 8164
$ws.Cells.Item(85, 12).Value = 8164
$ws.Cells.Item(85, 14).Value = -10660

$ws.Cells.Item(126, 8).Value = 8317.454
$ws.Cells.Item(126, 9).Value = 2566.8
$ws.Cells.Item(126, 11).Value = 7700.400000000001
$ws.Cells.Item(126, 13).Value = -5230.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 5960.6924
$ws.Cells.Item(62, 9).Value = 5279
$ws.Cells.Item(62, 10).Value = 6386.75
$ws.Cells.Item(62, 11).Value = 5279
$ws.Cells.Item(62, 12).Value = 6386.75
$ws.Cells.Item(62, 13).Value = -4655
$ws.Cells.Item(62, 14).Value = -7634.75

$ws.Cells.Item(65, 8).Value = 5960.6924
$ws.Cells.Item(65, 9).Value = 5279
$ws.Cells.Item(65, 10).Value = 6386.75
$ws.Cells.Item(65, 11).Value = 26395
$ws.Cells.Item(65, 12).Value = 31933.75
$ws.Cells.Item(65, 13).Value = -23275
$ws.Cells.Item(65, 14).Value = -38173.75

$ws.Cells.Item(104, 8).Value = 31551.334
$ws.Cells.Item(104, 10).Value = 31551.334
$ws.Cells.Item(104, 12).Value = 31551.334
$ws.Cells.Item(104, 14).Value = -38539.334

$ws.Cells.Item(105, 8).Value = 59500
$ws.Cells.Item(105, 10).Value = 59500
$ws.Cells.Item(105, 12).Value = 59500
$ws.Cells.Item(105, 14).Value = -66488

$ws.Cells.Item(107, 8).Value = 66667028
$ws.Cells.Item(107, 10).Value = 500000060
$ws.Cells.Item(107, 12).Value = 1500000180
$ws.Cells.Item(107, 14).Value = -1500004020

$ws.Cells.Item(126, 8).Value = 2659
$ws.Cells.Item(126, 9).Value = 2118.8462
$ws.Cells.Item(126, 11).Value = 6356.5386
$ws.Cells.Item(126, 13).Value = -3886.5386

$ws.Cells.Item(132, 8).Value = 25455.928
$ws.Cells.Item(132, 9).Value = 29268.723
$ws.Cells.Item(132, 11).Value = 87806.16900000001
$ws.Cells.Item(132, 13).Value = -85276.16900000001

Write-Host "Updated 44 leve rows across 8 job sheets."
